$wb = $excel.ActiveWorkbook

# "Overview" sheet: row 7 (f5fb7300-6c15-47f9-8ae9-084502f3aaa9.md) -
# the "Latest HO Xliff Generate Date" column (G) gets a fresh handback
# generation timestamp.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-28 14:41:32"

# "zh-cn" sheet: row 7 (same file) - "Latest Handback DateTime" column (H)
# is refreshed to the new handoff-report generation time.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-28 14:41:27"

# "de-de" sheet: row 7 (same file) - "Latest Handback DateTime" column (H)
# is refreshed to the new handoff-report generation time.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-28 14:41:32"
